# Update "Horarios" workbook - Linea 141 - scrape refresh 03:45:24
$wb = $excel.ActiveWorkbook

$newTime = "03:45:24"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 8"

$rows1 = @(
    @($newTime, "03:46", "14_ABASTO", 1, "LP1912"),
    @($newTime, "04:01", "81_EL PELIGRO", 16, "LP1912"),
    @($newTime, "04:46", "215A_EL PATO", 61, "LP1912"),
    @($newTime, "04:53", "11_ETCHEVERRY", 68, "LP1912"),
    @($newTime, "05:16", "17_ROMERO", 91, "LP1912"),
    @($newTime, "05:22", "23_HERNANDEZ", 97, "LP1912"),
    @($newTime, "05:34", "215B_EL PATO", 109, "LP1912"),
    @($newTime, "05:36", "14_ABASTO", 111, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"

$rows2 = @(
    @($newTime, "04:46", "215A_EL PATO", 61, "LP1912"),
    @($newTime, "05:34", "215B_EL PATO", 109, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws3.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws3.Cells.Item(5, 3).Value = "Linea"
$ws3.Cells.Item(5, 4).Value = "Minutos"
$ws3.Cells.Item(5, 5).Value = "Parada"

$ws3.Cells.Item(6, 1).Value = $newTime
$ws3.Cells.Item(6, 2).Value = "05:44"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 119
$ws3.Cells.Item(6, 5).Value = "L6173"

Write-Host "Horarios actualizados Linea 141"
